# Update countries & provincias Spain
# Applies the data refresh described by the commit: some country rows swap
# their country label (because the underlying country list was re-sorted)
# and a number of rows receive updated case statistics. The "last updated"
# timestamp footer is also bumped.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Country name corrections (rows whose label moved because of the
#     upstream re-sort of the source country list) ---
$ws.Range("A53").Value  = "Luxemburgo"
$ws.Range("A54").Value  = "Egipto"

$ws.Range("A90").Value  = "Guinea"
$ws.Range("A91").Value  = "Republica de Chipre"
$ws.Range("A92").Value  = "Letonia"

$ws.Range("A113").Value = "Sri Lanka"
$ws.Range("A114").Value = "Guatemala"

$ws.Range("A116").Value = "Somalia"
$ws.Range("A117").Value = "Mayotte"
$ws.Range("A118").Value = "Kenia"
$ws.Range("A119").Value = "Montenegro"
$ws.Range("A120").Value = "Mali"
$ws.Range("A121").Value = "Isla de Man"
$ws.Range("A122").Value = "Venezuela"

# --- Updated statistics (Casos totales, Nuevos casos, Casos activos,
#     Recuperados, Casos criticos, Muertes hoy, Muertes) ---

# Estados Unidos
$ws.Range("B4").Value  = 852717
$ws.Range("C4").Value  = 4000
$ws.Range("D4").Value  = 84221
$ws.Range("E4").Value  = 720196
$ws.Range("G4").Value  = 641
$ws.Range("H4").Value  = 48300

# Alemania
$ws.Range("B8").Value  = 151195
$ws.Range("C8").Value  = 547
$ws.Range("E8").Value  = 42541

# Singapur
$ws.Range("D29").Value = 924
$ws.Range("E29").Value = 10242
$ws.Range("F29").Value = 26

# Luxemburgo (row 53, formerly Egipto)
$ws.Range("B53").Value = 3665
$ws.Range("C53").Value = 11
$ws.Range("D53").Value = 711
$ws.Range("E53").Value = 2871
$ws.Range("F53").Value = 27
$ws.Range("G53").Value = 3
$ws.Range("H53").Value = 83

# Egipto (row 54, formerly Luxemburgo)
$ws.Range("B54").Value = 3659
$ws.Range("D54").Value = 935
$ws.Range("E54").Value = 2448
$ws.Range("F54").Value = 0
$ws.Range("H54").Value = 276

# Guinea (row 90, formerly Republica de Chipre)
$ws.Range("B90").Value = 862
$ws.Range("C90").Value = 101
$ws.Range("D90").Value = 170
$ws.Range("E90").Value = 686
$ws.Range("F90").Value = 0
$ws.Range("H90").Value = 6

# Republica de Chipre (row 91, formerly Letonia)
$ws.Range("B91").Value = 795
$ws.Range("C91").Value = 5
$ws.Range("D91").Value = 98
$ws.Range("E91").Value = 684
$ws.Range("F91").Value = 15
$ws.Range("H91").Value = 13

# Letonia (row 92, formerly Guinea)
$ws.Range("B92").Value = 778
$ws.Range("C92").Value = 17
$ws.Range("D92").Value = 133
$ws.Range("E92").Value = 634
$ws.Range("F92").Value = 6
$ws.Range("H92").Value = 11

# Sri Lanka (row 113, formerly Guatemala)
$ws.Range("B113").Value = 368
$ws.Range("C113").Value = 38
$ws.Range("D113").Value = 107
$ws.Range("E113").Value = 254
$ws.Range("F113").Value = 2
$ws.Range("G113").Value = 0
$ws.Range("H113").Value = 7

# Guatemala (row 114, formerly Sri Lanka)
$ws.Range("B114").Value = 342
$ws.Range("C114").Value = 26
$ws.Range("D114").Value = 25
$ws.Range("E114").Value = 307
$ws.Range("F114").Value = 3
$ws.Range("G114").Value = 2
$ws.Range("H114").Value = 10

# Somalia (row 116, formerly Mayotte)
$ws.Range("B116").Value = 328
$ws.Range("C116").Value = 42
$ws.Range("D116").Value = 8
$ws.Range("E116").Value = 304
$ws.Range("F116").Value = 2
$ws.Range("G116").Value = 8
$ws.Range("H116").Value = 16

# Mayotte (row 117, formerly Kenia)
$ws.Range("B117").Value = 326
$ws.Range("C117").Value = 0
$ws.Range("D117").Value = 125
$ws.Range("E117").Value = 197
$ws.Range("F117").Value = 4
$ws.Range("H117").Value = 4

# Kenia (row 118, formerly Montenegro)
$ws.Range("B118").Value = 320
$ws.Range("C118").Value = 17
$ws.Range("D118").Value = 89
$ws.Range("E118").Value = 217
$ws.Range("F118").Value = 2
$ws.Range("H118").Value = 14

# Montenegro (row 119, formerly Mali)
$ws.Range("B119").Value = 316
$ws.Range("C119").Value = 1
$ws.Range("D119").Value = 123
$ws.Range("E119").Value = 188
$ws.Range("F119").Value = 7
$ws.Range("G119").Value = 0
$ws.Range("H119").Value = 5

# Mali (row 120, formerly Isla de Man)
$ws.Range("B120").Value = 309
$ws.Range("C120").Value = 16
$ws.Range("D120").Value = 77
$ws.Range("E120").Value = 211
$ws.Range("F120").Value = 0
$ws.Range("G120").Value = 4
$ws.Range("H120").Value = 21

# Isla de Man (row 121, formerly Venezuela)
$ws.Range("B121").Value = 307
$ws.Range("D121").Value = 221
$ws.Range("E121").Value = 70
$ws.Range("F121").Value = 20
$ws.Range("G121").Value = 1
$ws.Range("H121").Value = 16

# Venezuela (row 122, formerly Somalia)
$ws.Range("B122").Value = 298
$ws.Range("D122").Value = 122
$ws.Range("E122").Value = 166
$ws.Range("F122").Value = 4
$ws.Range("H122").Value = 10

# --- Footer timestamp ---
$ws.Range("A1").Value = "Datos actualizados a 23 de Abril de 2020 a las 18:22"
